$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (commissioner) now points at the approver that used to live on the
# last row (old row 19, "commissioner2" / Ravindra Babu ~ ADM_Commissioner_2).
$ws.Range("D5").Value = "Ravindra Babu ~ ADM_Commissioner_2"

# Old row 14 was a duplicate of row 6 (commissioner1) - remove it, shifting
# rows 15-19 up by one.
$ws.Rows(14).Delete()

# Old row 19 (commissioner2, now at row 18 after the previous delete) is no
# longer needed since its data moved into row 5 above - remove it too.
$ws.Rows(18).Delete()

# Restore the active cell/selection left by the author after the edits.
$ws.Range("B22").Select()
